# "replaced "piece" observations in dataset with English titles and
# adjusted code"
#
# The "all_docs_lowercase" sheet's column P ("piece") held German-language
# values ("Massnahme" / "Mutter") used as observation labels. Replace them
# with their English equivalents ("measures" / "mother").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_docs_lowercase")

# Rows 2-51 and 133-134 held "Massnahme" ("measure") in column P;
# replace with the English title "measures".
$ws.Range("P2:P51").Value = "measures"
$ws.Range("P133:P134").Value = "measures"

# Rows 52-132 held "Mutter" ("mother") in column P; replace with the
# English title "mother".
$ws.Range("P52:P132").Value = "mother"

# Bring the sheet to the front and scroll/select to match the author's
# final view (scrolled right to column I, active cell M15).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M15").Select()
